$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.865.66'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '1.841.01'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'231.58"
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'39.92"
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").Value = '2.108.98'
$ws.Range("E12").Value = '  +1.93%  '
$ws.Range("D13").Value = "'11.64"
$ws.Range("E13").Value = '  +5.56%  '
$ws.Range("D14").Value = '1.847.56'
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '34.892.44'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = "'69.89"
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = "'240.43"
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("E21").Value = '  +2.34%  '
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = "'2.27"
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  -1.04%  '
$ws.Range("D27").Value = "'17.48"
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("E29").Value = '  -2.79%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").Value = "'3.95"
$ws.Range("E32").Value = '  -4.61%  '
$ws.Range("D33").Value = "'3.97"
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("E34").Value = '  +8.12%  '
$ws.Range("E35").Value = '  +7.23%  '
$ws.Range("E36").Value = '  +13.04%  '
$ws.Range("D37").Value = "'0.697"
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("E38").Value = '  +7.75%  '
$ws.Range("D39").Value = "'90.32"
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("D40").Value = '1.349.56'
$ws.Range("E40").Value = '  +2.94%  '
$ws.Range("E41").Value = '  +0.47%  '
$ws.Range("D42").Value = "'14.85"
$ws.Range("E42").Value = '  +2.97%  '
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").Value = "'2.41"
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("E47").Value = '  +2.16%  '
$ws.Range("D48").Value = '2.022.46'
$ws.Range("E48").Value = '  +1.85%  '
$ws.Range("E49").Value = '  +23.34%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("E51").Value = '  +0.63%  '
